# The sheet currently has two "merge & center" blocks in the first two
# columns: A2:A4 (a date) and B2:B4 (a shift name, "Morning"). Each block
# only carries its value in the top-left cell (A2 / B2); the other rows in
# the block are blank and simply show the merged value.
#
# This edit un-merges both blocks and fills the value that used to be
# shown (via the merge) into every row of the block, so A2:A4 each hold
# the date and B2:B4 each hold "Morning" as their own, independent cell
# values - matching rows 5 and 6 below, which were never merged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the values that are currently visible (held by the top-left
# cell of each merged block) before unmerging.
$dateVal = $ws.Range("A2").Value()
$shiftVal = $ws.Range("B2").Value()

# Unmerge the two blocks.
$ws.Range("A2:A4").UnMerge()
$ws.Range("B2:B4").UnMerge()

# Fill the previously-merged value into every row of each block.
$ws.Range("A3").Value = $dateVal
$ws.Range("A4").Value = $dateVal
$ws.Range("B3").Value = $shiftVal
$ws.Range("B4").Value = $shiftVal

# The merge-and-center formatting had centered the date column
# horizontally; drop that now that the cells stand alone (matching the
# plain, vertical-center-only look used by the date cells further down
# the column).
$ws.Range("A2:A6").HorizontalAlignment = 1

# Leave the same cell selected/active as when the workbook was last saved.
$ws.Range("D17").Select()
